$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Cell values in column D that would otherwise be auto-coerced to numbers
# by Excel (because they parse as plain numeric literals) need to be forced
# to remain text, matching the original inlineStr/text storage for that column.

$ws.Range('D2').Value = '29.244.53'
$ws.Range('E2').Value = '  -0.07%  '
$ws.Range('D3').Value = '1.842.89'
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.9991'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '241.02'
$ws.Range('E5').Value = '  -1.00%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.6698'
$ws.Range('E6').Value = '  -2.50%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.9997'
$ws.Range('E7').Value = '  -0.10%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.07462'
$ws.Range('E8').Value = '  -0.77%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.2964'
$ws.Range('E9').Value = '  -2.17%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '22.77'
$ws.Range('E10').Value = '  -2.12%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07721'
$ws.Range('E11').Value = '  +0.38%  '
$ws.Range('E12').Value = '  -1.27%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.6785'
$ws.Range('E13').Value = '  -1.17%  '
$ws.Range('D14').Value = '1.776.35'
$ws.Range('E14').Value = '  -3.47%  '
$ws.Range('E15').Value = '  -2.50%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '6.189'
$ws.Range('E16').Value = '  -1.40%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.000008329'
$ws.Range('E17').Value = '  +1.38%  '
$ws.Range('D18').Value = '29.041.00'
$ws.Range('E18').Value = '  -0.74%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '228.83'
$ws.Range('E19').Value = '  -1.57%  '
$ws.Range('E20').Value = '  -0.41%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '0.9992'
$ws.Range('E21').Value = '  -0.11%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '7.197'
$ws.Range('E22').Value = '  -3.59%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.9993'
$ws.Range('E23').Value = '  -0.08%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '159.94'
$ws.Range('E24').Value = '  -0.01%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '8.703'
$ws.Range('E25').Value = '  -1.52%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.1401'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '18.04'
$ws.Range('E27').Value = '  -0.45%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.512'
$ws.Range('E28').Value = '  -0.52%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '4.193'
$ws.Range('E29').Value = '  -1.64%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '4.090'
$ws.Range('E30').Value = '  -1.59%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.190'
$ws.Range('E31').Value = '  -1.61%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.05336'
$ws.Range('E32').Value = '  +3.55%  '
$ws.Range('B33').Value = 'LidoDAOToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.889'
$ws.Range('E33').Value = '  +2.63%  '
$ws.Range('B34').Value = 'ImmutableX'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.7601'
$ws.Range('E34').Value = '  -1.98%  '
$ws.Range('E35').Value = '  +0.40%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.675'
$ws.Range('E36').Value = '  +0.01%  '
$ws.Range('D37').Value = '1.334.99'
$ws.Range('E37').Value = '  +2.27%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.01808'
$ws.Range('E38').Value = '  -1.91%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.738'
$ws.Range('E39').Value = '  +1.32%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.9232'
$ws.Range('E40').Value = '  -2.15%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '5.952'
$ws.Range('E41').Value = '  +3.07%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.001'
$ws.Range('E42').Value = '  +0.10%  '
$ws.Range('B43').Value = 'Quant'
$ws.Range('C43').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '103.41'
$ws.Range('E43').Value = '  -1.59%  '
$ws.Range('B44').Value = 'XinFinNetwork'
$ws.Range('C44').Value = 'https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.08154'
$ws.Range('E44').Value = '  +15.10%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.00000000124'
$ws.Range('E45').Value = '  +2.00%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.5158'
$ws.Range('E46').Value = '  -0.78%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.768'
$ws.Range('E47').Value = '  -0.33%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '63.74'
$ws.Range('E48').Value = '  -1.38%  '
$ws.Range('D49').Value = '1.919.09'
$ws.Range('E49').Value = '  -3.46%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '9.227'
$ws.Range('E50').Value = '  -4.33%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.05931'
$ws.Range('E51').Value = '  +0.16%  '
